# Add a new "percent_of_control" column (D) to Sheet3, computing each
# treatment's average Chl_a as a percentage of the control (column B
# divided by the control value in I1), mirroring the existing
# "percent_change_relative_to_control" column already in C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("D1").Value = "percent_of_control"

# New formulas: (treatment_avg / control_avg) * 100
$ws.Range("D2").Formula = '=(B2/$I$1)*100'
$ws.Range("D3").Formula = '=(B3/$I$1)*100'
$ws.Range("D4").Formula = '=(B4/$I$1)*100'
$ws.Range("D5").Formula = '=(B5/$I$1)*100'
$ws.Range("D6").Formula = '=(B6/$I$1)*100'

# Match the last selected cell recorded in the workbook
$ws.Range("D6").Select()
